# Generate Report for Handoff
# For files 22ed719f, 6f54e31e, d2e2485e, e35c5b91 (rows 4-7 in both the
# zh-cn and de-de status tables), a new handoff round completed:
#   - Priority flips from "low" to "ht"
#   - Latest Handoff Datetime is refreshed to the new generation time
# The zh-cn sheet's handoff time moves to 2016-09-03 02:34:26, the de-de
# sheet's (and therefore the Overview rollup's) moves to 2016-09-03 02:34:30.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $ws_zhcn.Cells.Item($r, 5).Value = "ht"
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-09-03 02:34:26"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $ws_dede.Cells.Item($r, 5).Value = "ht"
    $ws_dede.Cells.Item($r, 8).Value = "2016-09-03 02:34:30"
}

# The Overview rollup sheet mirrors the de-de "Latest Handoff Datetime"
# (it shared the same underlying string as de-de!H4:H7 before the edit),
# so keep it in sync too.
$ws_overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $ws_overview.Cells.Item($r, 7).Value = "2016-09-03 02:34:30"
}
